$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.455765
$ws.Range("H2").Value = 13.367295
$ws.Range("I2").Value = 0.1558824083674925
$ws.Range("J2").Value = 0.167793131187596
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.5297143333333333
$ws.Range("N2").Value = 1.589143
$ws.Range("O2").Value = 0.03463898003508423
$ws.Range("P2").Value = 0.0369789110027492
$ws.Range("Q2").Value = 2.360282586464999
$ws.Range("R2").Value = 21.242543278185
$ws.Range("S2").Value = 0.005399607631262421
$ws.Range("T2").Value = 0.006204807265058735
$ws.Range("G3").Value = 4.455765
$ws.Range("H3").Value = 13.367295
$ws.Range("I3").Value = 0.1558824083674925
$ws.Range("J3").Value = 0.167793131187596
$ws.Range("O3").Value = 0.7755286881671239
$ws.Range("P3").Value = 0.8279171705045608
$ws.Range("Q3").Value = 52.84413271207499
$ws.Range("R3").Value = 475.597194408675
$ws.Range("S3").Value = 0.1208912796695734
$ws.Range("T3").Value = 0.1389188144029351
$ws.Range("G4").Value = 4.455765
$ws.Range("H4").Value = 13.367295
$ws.Range("I4").Value = 0.1558824083674925
$ws.Range("J4").Value = 0.167793131187596
$ws.Range("M4").Value = 2.9029985
$ws.Range("N4").Value = 5.805997
$ws.Range("O4").Value = 0.1898323317977919
$ws.Range("P4").Value = 0.13510391849269
$ws.Range("Q4").Value = 12.9350791113525
$ws.Range("R4").Value = 77.61047466811499
$ws.Range("S4").Value = 0.02959152106665673
$ws.Range("T4").Value = 0.02266950951960221
$ws.Range("I5").Value = 0.1858758098371279
$ws.Range("J5").Value = 0.2000782799754709
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.5297143333333333
$ws.Range("N5").Value = 1.589143
$ws.Range("O5").Value = 0.03463898003508423
$ws.Range("P5").Value = 0.0369789110027492
$ws.Range("Q5").Value = 2.814425577576222
$ws.Range("R5").Value = 25.329830198186
$ws.Range("S5").Value = 0.006438548465953387
$ws.Range("T5").Value = 0.007398676908796076
$ws.Range("I6").Value = 0.1858758098371279
$ws.Range("J6").Value = 0.2000782799754709
$ws.Range("O6").Value = 0.7755286881671239
$ws.Range("P6").Value = 0.8279171705045608
$ws.Range("S6").Value = 0.1441520229649896
$ws.Range("T6").Value = 0.1656482434367112
$ws.Range("I7").Value = 0.1858758098371279
$ws.Range("J7").Value = 0.2000782799754709
$ws.Range("M7").Value = 2.9029985
$ws.Range("N7").Value = 5.805997
$ws.Range("O7").Value = 0.1898323317977919
$ws.Range("P7").Value = 0.13510391849269
$ws.Range("Q7").Value = 15.42392326568233
$ws.Range("R7").Value = 92.543539594094
$ws.Range("S7").Value = 0.03528523840618494
$ws.Range("T7").Value = 0.02703135962996363
$ws.Range("G8").Value = 5.633732333333334
$ws.Range("H8").Value = 16.901197
$ws.Range("I8").Value = 0.1970929266282699
$ws.Range("J8").Value = 0.2121524785267629
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.5297143333333333
$ws.Range("N8").Value = 1.589143
$ws.Range("O8").Value = 0.03463898003508423
$ws.Range("P8").Value = 0.0369789110027492
$ws.Range("Q8").Value = 2.984268767130111
$ws.Range("R8").Value = 26.858418904171
$ws.Range("S8").Value = 0.006827097950532964
$ws.Range("T8").Value = 0.007845167622453825
$ws.Range("G9").Value = 5.633732333333334
$ws.Range("H9").Value = 16.901197
$ws.Range("I9").Value = 0.1970929266282699
$ws.Range("J9").Value = 0.2121524785267629
$ws.Range("O9").Value = 0.7755286881671239
$ws.Range("P9").Value = 0.8279171705045608
$ws.Range("Q9").Value = 66.81449741783389
$ws.Range("R9").Value = 601.330476760505
$ws.Range("S9").Value = 0.1528512188350414
$ws.Range("T9").Value = 0.1756446797374071
$ws.Range("G10").Value = 5.633732333333334
$ws.Range("H10").Value = 16.901197
$ws.Range("I10").Value = 0.1970929266282699
$ws.Range("J10").Value = 0.2121524785267629
$ws.Range("M10").Value = 2.9029985
$ws.Range("N10").Value = 5.805997
$ws.Range("O10").Value = 0.1898323317977919
$ws.Range("P10").Value = 0.13510391849269
$ws.Range("Q10").Value = 16.35471651306817
$ws.Range("R10").Value = 98.12829907840899
$ws.Range("S10").Value = 0.0374146098426956
$ws.Range("T10").Value = 0.02866263116690194
$ws.Range("G11").Value = 6.087099
$ws.Range("H11").Value = 12.174198
$ws.Range("I11").Value = 0.2129537020222914
$ws.Range("J11").Value = 0.15281676675182
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.5297143333333333
$ws.Range("N11").Value = 1.589143
$ws.Range("O11").Value = 0.03463898003508423
$ws.Range("P11").Value = 0.0369789110027492
$ws.Range("Q11").Value = 3.224423588719
$ws.Range("R11").Value = 19.346541532314
$ws.Range("S11").Value = 0.007376499032747431
$ws.Range("T11").Value = 0.005650997617443435
$ws.Range("G12").Value = 6.087099
$ws.Range("H12").Value = 12.174198
$ws.Range("I12").Value = 0.2129537020222914
$ws.Range("J12").Value = 0.15281676675182
$ws.Range("O12").Value = 0.7755286881671239
$ws.Range("P12").Value = 0.8279171705045608
$ws.Range("Q12").Value = 72.19129989744501
$ws.Range("R12").Value = 433.1477993846701
$ws.Range("S12").Value = 0.1651517051696803
$ws.Range("T12").Value = 0.1265196251348223
$ws.Range("G13").Value = 6.087099
$ws.Range("H13").Value = 12.174198
$ws.Range("I13").Value = 0.2129537020222914
$ws.Range("J13").Value = 0.15281676675182
$ws.Range("M13").Value = 2.9029985
$ws.Range("N13").Value = 5.805997
$ws.Range("O13").Value = 0.1898323317977919
$ws.Range("P13").Value = 0.13510391849269
$ws.Range("Q13").Value = 17.6708392663515
$ws.Range("R13").Value = 70.683357065406
$ws.Range("S13").Value = 0.04042549781986374
$ws.Range("T13").Value = 0.02064614399955431
$ws.Range("G14").Value = 7.094445666666666
$ws.Range("H14").Value = 21.283337
$ws.Range("I14").Value = 0.2481951531448182
$ws.Range("J14").Value = 0.2671593435583502
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.5297143333333333
$ws.Range("N14").Value = 1.589143
$ws.Range("O14").Value = 0.03463898003508423
$ws.Range("P14").Value = 0.0369789110027492
$ws.Range("Q14").Value = 3.758029556687888
$ws.Range("R14").Value = 33.822266010191
$ws.Range("S14").Value = 0.008597226954588033
$ws.Range("T14").Value = 0.009879261588997129
$ws.Range("G15").Value = 7.094445666666666
$ws.Range("H15").Value = 21.283337
$ws.Range("I15").Value = 0.2481951531448182
$ws.Range("J15").Value = 0.2671593435583502
$ws.Range("O15").Value = 0.7755286881671239
$ws.Range("P15").Value = 0.8279171705045608
$ws.Range("Q15").Value = 84.1381509859561
$ws.Range("R15").Value = 757.2433588736051
$ws.Range("S15").Value = 0.1924824615278393
$ws.Range("T15").Value = 0.2211858077926852
$ws.Range("G16").Value = 7.094445666666666
$ws.Range("H16").Value = 21.283337
$ws.Range("I16").Value = 0.2481951531448182
$ws.Range("J16").Value = 0.2671593435583502
$ws.Range("M16").Value = 2.9029985
$ws.Range("N16").Value = 5.805997
$ws.Range("O16").Value = 0.1898323317977919
$ws.Range("P16").Value = 0.13510391849269
$ws.Range("Q16").Value = 20.59516512866483
$ws.Range("R16").Value = 123.570990771989
$ws.Range("S16").Value = 0.04711546466239092
$ws.Range("T16").Value = 0.03609427417666791
